$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits at the very end of the document
# (after the last existing entry, "- Horváth Bence János"). The new
# content being appended moves that bookmark so it ends up inside the
# third newly-added entry (the one about collecting/generating images),
# right after "- Kovács Milán". Remove the old bookmark first so we can
# re-create it in its new location as part of the inserted XML below.
if ($d.Bookmarks.Exists("_GoBack")) {
  $d.Bookmarks.Item("_GoBack").Delete()
}

# Collapse a range to the very end of the document body (just before
# the final sectPr) so the new paragraphs are appended after the last
# existing project-log entry.
$r = $d.Content
$r.Collapse(0)

$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @"
<w:p $ns>
  <w:pPr>
    <w:keepLines/>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">2025. március 5. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>Megbeszélés a frontend implementációjáról</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>- Balogh Attila, Horváth Bence János, Kovács Milán</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:keepLines/>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>2025. március 8.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>Főoldal implementálása frontend oldalon</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>- Balogh Attila</w:t>
  </w:r>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:keepLines/>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>2025. március 8.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>Weboldalhoz szükséges képek összegyűjtése és generálása</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>- Kovács Milán</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $ns>
  <w:pPr>
    <w:keepLines/>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:t>2025. március 11.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>Mellékoldalak implementációjának kezdete frontend oldalon</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="24"/>
    </w:rPr>
    <w:br/>
    <w:t>- Balogh Attila</w:t>
  </w:r>
</w:p>
"@

$r.InsertXML($xml)
